# Applies the data reshuffle described by the diff: the SymbolVals rows
# (A2:F25) were regenerated to include maxVal & minVal, which changed the
# resulting row order / values for the "symbol" data table on the active
# sheet. Row 1 (headers) and row 26 (totals) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(901,  16, 15, 45, 60, 60),
    @(301,  6,  45, 30, 60, 45),
    @(1202, 2,  10, 10, 10, 10),
    @(1203, 3,  15, 15, 15, 15),
    @(801,  3,  67, 65, 52, 45),
    @(902,  1,  0,  0,  0,  0),
    @(401,  9,  48, 67, 75, 45),
    @(1201, 2,  10, 10, 10, 10),
    @(601,  9,  60, 67, 60, 42),
    @(201,  9,  30, 15, 45, 30),
    @(101,  9,  30, 15, 60, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(501,  9,  52, 30, 75, 45),
    @(701,  3,  90, 45, 97, 15),
    @(1101, 0,  15, 30, 30, 0),
    @(1,    0,  2,  2,  2,  2),
    @(2,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3),
    @(502,  0,  4,  0,  0,  0),
    @(802,  0,  4,  5,  4,  0),
    @(602,  0,  0,  4,  0,  9),
    @(402,  0,  0,  4,  0,  0),
    @(702,  0,  0,  0,  4,  0),
    @(1002, 0,  0,  0,  0,  9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowVals = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowVals[$col - 1]
    }
}
